$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the time-worked values for row 28 (11-11-2020)
$ws.Range("B28").Value = 10
$ws.Range("C28").Value = 13

# Add descriptions for rows 28 and 29
$ws.Range("E28").Value = "Trying to fix some issues with displaying the parsed JSON data"
$ws.Range("E29").Value = "Maalisuora Event"

# Move the active selection to E29 (as saved in the workbook view)
$ws.Range("E29").Select()
